$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 6885.524
$ws.Range("J17").Value = 7717.4116
$ws.Range("L17").Value = 23152.2348
$ws.Range("N17").Value = -23488.2348
$ws.Range("H40").Value = 6395.364
$ws.Range("J40").Value = 6424.8335
$ws.Range("L40").Value = 6424.8335
$ws.Range("N40").Value = -6774.8335
$ws.Range("H74").Value = 6734.913
$ws.Range("I74").Value = 5361.0557
$ws.Range("K74").Value = 5361.0557
$ws.Range("M74").Value = -4425.0557
$ws.Range("H77").Value = 6734.913
$ws.Range("I77").Value = 5361.0557
$ws.Range("K77").Value = 26805.2785
$ws.Range("M77").Value = -22125.2785
$ws.Range("H98").Value = 382620.44
$ws.Range("I98").Value = 1532.8462
$ws.Range("J98").Value = 2034000
$ws.Range("K98").Value = 1532.8462
$ws.Range("L98").Value = 2034000
$ws.Range("M98").Value = -34.84619999999995
$ws.Range("N98").Value = -2036996
$ws.Range("H101").Value = 879.6316
$ws.Range("J101").Value = 1145.625
$ws.Range("L101").Value = 3436.875
$ws.Range("N101").Value = -6680.875
$ws.Range("H112").Value = 1565.2963
$ws.Range("I112").Value = 1026.6666
$ws.Range("J112").Value = 1632.625
$ws.Range("K112").Value = 3079.9998
$ws.Range("L112").Value = 4897.875
$ws.Range("M112").Value = -1971.9998
$ws.Range("N112").Value = -7113.875
$ws.Range("H122").Value = 382620.44
$ws.Range("I122").Value = 1532.8462
$ws.Range("J122").Value = 2034000
$ws.Range("K122").Value = 4598.5386
$ws.Range("L122").Value = 6102000
$ws.Range("M122").Value = -2148.5386
$ws.Range("N122").Value = -6106900
$ws.Range("H125").Value = 2543.3333
$ws.Range("I125").Value = 1485.5
$ws.Range("J125").Value = 3389.6
$ws.Range("K125").Value = 13369.5
$ws.Range("L125").Value = 30506.4
$ws.Range("M125").Value = -10909.5
$ws.Range("N125").Value = -35426.39999999999
$ws.Range("H138").Value = 2599.7454
$ws.Range("I138").Value = 1076.7812
$ws.Range("K138").Value = 3230.3436
$ws.Range("M138").Value = 1909.6564
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 70001
$ws.Range("J7").Value = 70001
$ws.Range("L7").Value = 70001
$ws.Range("N7").Value = -70229
$ws.Range("H32").Value = 3158.4033
$ws.Range("I32").Value = 2777.034
$ws.Range("K32").Value = 2777.034
$ws.Range("M32").Value = -2490.034
$ws.Range("H45").Value = 1852.3125
$ws.Range("I45").Value = 1857.5333
$ws.Range("K45").Value = 1857.5333
$ws.Range("M45").Value = -1480.5333
$ws.Range("H61").Value = 6304.7856
$ws.Range("I61").Value = 2854.4167
$ws.Range("J61").Value = 27007
$ws.Range("K61").Value = 2854.4167
$ws.Range("L61").Value = 27007
$ws.Range("M61").Value = -2642.4167
$ws.Range("N61").Value = -27431
$ws.Range("H88").Value = 4616.778
$ws.Range("I88").Value = 3805.5715
$ws.Range("K88").Value = 3805.5715
$ws.Range("M88").Value = -3399.5715
$ws.Range("H91").Value = 4616.778
$ws.Range("I91").Value = 3805.5715
$ws.Range("K91").Value = 3805.5715
$ws.Range("M91").Value = -2401.5715
$ws.Range("H136").Value = 6304.7856
$ws.Range("I136").Value = 2854.4167
$ws.Range("J136").Value = 27007
$ws.Range("K136").Value = 8563.250100000001
$ws.Range("L136").Value = 81021
$ws.Range("M136").Value = -6013.250100000001
$ws.Range("N136").Value = -86121
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 101999
$ws.Range("J55").Value = 101999
$ws.Range("L55").Value = 101999
$ws.Range("N55").Value = -102545
$ws.Range("H134").Value = 3107.4285
$ws.Range("I134").Value = 3107.4285
$ws.Range("K134").Value = 9322.2855
$ws.Range("M134").Value = -6787.2855
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3051.4546
$ws.Range("I99").Value = 2544.3333
$ws.Range("K99").Value = 2544.3333
$ws.Range("M99").Value = -1046.3333
$ws.Range("H126").Value = 3051.4546
$ws.Range("I126").Value = 2544.3333
$ws.Range("K126").Value = 7632.999899999999
$ws.Range("M126").Value = -5162.999899999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1970.6154
$ws.Range("I34").Value = 1827.875
$ws.Range("J34").Value = 2199
$ws.Range("K34").Value = 5483.625
$ws.Range("L34").Value = 6597
$ws.Range("M34").Value = -5399.625
$ws.Range("N34").Value = -6765
$ws.Range("H131").Value = 6798567
$ws.Range("I131").Value = 13890799
$ws.Range("K131").Value = 41672397
$ws.Range("M131").Value = -41667357
$ws.Range("H132").Value = 3129
$ws.Range("I132").Value = 1082.8572
$ws.Range("J132").Value = 4230.769
$ws.Range("K132").Value = 9745.7148
$ws.Range("L132").Value = 38076.921
$ws.Range("M132").Value = -7215.7148
$ws.Range("N132").Value = -43136.921
$ws.Range("H136").Value = 2344.9092
$ws.Range("I136").Value = 2079.4
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 6238.200000000001
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -1138.200000000001
$ws.Range("N136").Value = -25200
$ws.Range("H139").Value = 3281.9473
$ws.Range("I139").Value = 2005.7693
$ws.Range("J139").Value = 6047
$ws.Range("K139").Value = 6017.3079
$ws.Range("L139").Value = 18141
$ws.Range("M139").Value = -877.3078999999998
$ws.Range("N139").Value = -28421
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 752.75
$ws.Range("J2").Value = 1707.6
$ws.Range("L2").Value = 1707.6
$ws.Range("N2").Value = -1933.6
$ws.Range("H122").Value = 4970.4326
$ws.Range("J122").Value = 7572.1665
$ws.Range("L122").Value = 22716.4995
$ws.Range("N122").Value = -27616.4995
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 27499.5
$ws.Range("I5").Value = 27499.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 27499.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -27386.5
$ws.Range("N5").ClearContents()
$ws.Range("H7").Value = 5751.2
$ws.Range("I7").Value = 4550.3057
$ws.Range("J7").Value = 10554.777
$ws.Range("K7").Value = 4550.3057
$ws.Range("L7").Value = 10554.777
$ws.Range("M7").Value = -4438.3057
$ws.Range("N7").Value = -10778.777
$ws.Range("H100").Value = 9288.75
$ws.Range("I100").Value = 2999
$ws.Range("K100").Value = 2999
$ws.Range("M100").Value = -2458
$ws.Range("H104").Value = 44059
$ws.Range("J104").Value = 44059
$ws.Range("L104").Value = 44059
$ws.Range("N104").Value = -51047
$ws.Range("H105").Value = 80000
$ws.Range("J105").Value = 80000
$ws.Range("L105").Value = 80000
$ws.Range("N105").Value = -86988
$ws.Range("H126").Value = 5751.2
$ws.Range("I126").Value = 4550.3057
$ws.Range("J126").Value = 10554.777
$ws.Range("K126").Value = 13650.9171
$ws.Range("L126").Value = 31664.331
$ws.Range("M126").Value = -11180.9171
$ws.Range("N126").Value = -36604.331
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H60").Value = 48000
$ws.Range("I60").Value = 48000
$ws.Range("K60").Value = 48000
$ws.Range("M60").Value = -47178
$ws.Range("H88").Value = 21094
$ws.Range("J88").Value = 21094
$ws.Range("L88").Value = 21094
$ws.Range("N88").Value = -21906
$ws.Range("H91").Value = 21094
$ws.Range("J91").Value = 21094
$ws.Range("L91").Value = 21094
$ws.Range("N91").Value = -23902
$ws.Range("H133").Value = 67500
$ws.Range("J133").Value = 67500
$ws.Range("L133").Value = 67500
$ws.Range("N133").Value = -77620
$ws.Range("H136").Value = 2094.17
$ws.Range("I136").Value = 1779.1372
$ws.Range("K136").Value = 5337.411599999999
$ws.Range("M136").Value = -2787.411599999999
